$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.801.79'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '3.087.53'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.31'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '610.01'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.09'
$ws.Range('E7').Value = '  -3.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.385'
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '3.083.50'
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.728'
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = '92.079.14'
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.86'
$ws.Range('E15').Value = '  -4.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.39'
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('D17').Value = '3.673.87'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '3.096.12'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  -2.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.54'
$ws.Range('E20').Value = '  -4.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.75'
$ws.Range('E21').Value = '  -5.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.24'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '441.76'
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000193'
$ws.Range('E24').Value = '  -5.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.70'
$ws.Range('E25').Value = '  -5.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '85.66'
$ws.Range('E26').Value = '  -3.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.53'
$ws.Range('E27').Value = '  -4.50%  '
$ws.Range('D28').Value = '3.257.72'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.130'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.229'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('E32').Value = '  -2.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.06'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.82'
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('E36').Value = '  -8.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.81'
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  -4.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '480.31'
$ws.Range('E40').Value = '  -6.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.87'
$ws.Range('E41').Value = '  +7.79%  '
$ws.Range('E42').Value = '  -5.89%  '
$ws.Range('E43').Value = '  -5.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.27'
$ws.Range('E44').Value = '  -6.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '160.86'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('E47').Value = '  -4.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.681'
$ws.Range('E48').Value = '  -5.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0331'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.81'
$ws.Range('E51').Value = '  -0.55%  '
